$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.347.67"
Set-TextValue $ws.Range("E2") "  +0.04%  "
Set-TextValue $ws.Range("D3") "3.686.56"
Set-TextValue $ws.Range("E3") "  +0.06%  "
Set-TextValue $ws.Range("E4") "  +0.10%  "
Set-TextValue $ws.Range("D5") "678.63"
Set-TextValue $ws.Range("E5") "  -0.83%  "
Set-TextValue $ws.Range("D6") "160.45"
Set-TextValue $ws.Range("E6") "  +0.27%  "
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.494"
Set-TextValue $ws.Range("E8") "  +0.28%  "
Set-TextValue $ws.Range("E9") "  +0.50%  "
Set-TextValue $ws.Range("D10") "7.15"
Set-TextValue $ws.Range("E10") "  -0.77%  "
Set-TextValue $ws.Range("D11") "0.438"
Set-TextValue $ws.Range("E11") "  +0.60%  "
Set-TextValue $ws.Range("E12") "  +0.17%  "
Set-TextValue $ws.Range("D13") "4.314.13"
Set-TextValue $ws.Range("E13") "  +0.23%  "
Set-TextValue $ws.Range("D14") "32.40"
Set-TextValue $ws.Range("E14") "  +0.19%  "
Set-TextValue $ws.Range("D15") "3.666.45"
Set-TextValue $ws.Range("E15") "  -0.61%  "
Set-TextValue $ws.Range("D16") "69.378.85"
Set-TextValue $ws.Range("E16") "  +0.07%  "
Set-TextValue $ws.Range("E17") "  +2.69%  "
Set-TextValue $ws.Range("D18") "15.99"
Set-TextValue $ws.Range("D19") "6.46"
Set-TextValue $ws.Range("E19") "  +0.57%  "
Set-TextValue $ws.Range("D20") "472.69"
Set-TextValue $ws.Range("E20") "  +0.04%  "
Set-TextValue $ws.Range("D21") "9.78"
Set-TextValue $ws.Range("E21") "  -0.99%  "
Set-TextValue $ws.Range("D22") "0.649"
Set-TextValue $ws.Range("E22") "  +0.62%  "
Set-TextValue $ws.Range("D23") "80.11"
Set-TextValue $ws.Range("E23") "  +0.75%  "
Set-TextValue $ws.Range("D24") "3.835.31"
Set-TextValue $ws.Range("E24") "  +0.17%  "
Set-TextValue $ws.Range("E25") "  -0.04%  "
Set-TextValue $ws.Range("E26") "  +0.46%  "
Set-TextValue $ws.Range("D27") "10.87"
Set-TextValue $ws.Range("E27") "  +0.07%  "
Set-TextValue $ws.Range("D28") "9.10"
Set-TextValue $ws.Range("E28") "  -0.53%  "
Set-TextValue $ws.Range("E29") "  +0.36%  "
Set-TextValue $ws.Range("D30") "1.73"
Set-TextValue $ws.Range("E30") "  -0.29%  "
Set-TextValue $ws.Range("D31") "2.01"
Set-TextValue $ws.Range("E31") "  -0.42%  "
Set-TextValue $ws.Range("D32") "6.57"
Set-TextValue $ws.Range("E32") "  -0.96%  "
Set-TextValue $ws.Range("E33") "  +0.27%  "
Set-TextValue $ws.Range("D34") "26.95"
Set-TextValue $ws.Range("E34") "  +1.28%  "
Set-TextValue $ws.Range("D35") "3.679.66"
Set-TextValue $ws.Range("E35") "  +0.56%  "
Set-TextValue $ws.Range("E36") "  +2.04%  "
Set-TextValue $ws.Range("D37") "8.43"
Set-TextValue $ws.Range("E37") "  +3.38%  "
Set-TextValue $ws.Range("D38") "6.19"
Set-TextValue $ws.Range("E38") "  +2.08%  "
Set-TextValue $ws.Range("D40") "2.25"
Set-TextValue $ws.Range("E40") "  -0.36%  "
Set-TextValue $ws.Range("E41") "  +0.09%  "
Set-TextValue $ws.Range("D42") "0.0903"
Set-TextValue $ws.Range("E42") "  +0.50%  "
Set-TextValue $ws.Range("D43") "168.47"
Set-TextValue $ws.Range("E43") "  +2.11%  "
Set-TextValue $ws.Range("D44") "0.940"
Set-TextValue $ws.Range("E44") "  +0.00%  "
Set-TextValue $ws.Range("D45") "46.63"
Set-TextValue $ws.Range("E45") "  -2.57%  "
Set-TextValue $ws.Range("D46") "28.22"
Set-TextValue $ws.Range("E46") "  +0.17%  "
Set-TextValue $ws.Range("D47") "2.71"
Set-TextValue $ws.Range("E47") "  +0.34%  "
Set-TextValue $ws.Range("D48") "0.000278"
Set-TextValue $ws.Range("E48") "  +2.48%  "
Set-TextValue $ws.Range("D49") "1.28"
Set-TextValue $ws.Range("E49") "  -1.38%  "
Set-TextValue $ws.Range("D50") "1.08"
Set-TextValue $ws.Range("E50") "  -1.69%  "
Set-TextValue $ws.Range("D51") "7.86"
Set-TextValue $ws.Range("E51") "  +0.34%  "
